$d = $word.ActiveDocument
$br = [char]11

# --- 1. Title (Heading1) ---
$d.Content.Find.Execute(
    "Review 154: Context is Environment, 26.09.2023",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Review 153: [Short] FOLEYGEN: VISUALLY-GUIDED AUDIO GENERATION, 24.09.2023",
    2) | Out-Null

# --- 2. Paper link (bold line) ---
$d.Content.Find.Execute(
    "Paper: https://arxiv.org/abs/2309.09888v2",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Paper: https://arxiv.org/abs/2309.10537v1",
    2) | Out-Null

# --- 3. Merge the two closing paragraphs ("כאשר מודל..." + "החוקרים מראים...")
#        into a single paragraph, then overwrite its text with the new review
#        body. Work from the bottom of the document upward so the indices of
#        paragraphs above stay valid for the later steps.
$pKasher = $d.Paragraphs.Item(11)
$pHokrim = $d.Paragraphs.Item(12)
$markStart = $pKasher.Range.End - 1
$markEnd = $pKasher.Range.End
$d.Range($markStart, $markEnd).Delete()

$pMerged2 = $d.Paragraphs.Item(11)
$s = $pMerged2.Range.Start
$e = $pMerged2.Range.End - 1
$newText2 = "בחרתי לסקור את המאמר הזה כי למדתי ממנו שיש משימה שלא הכרתי בראייה ממוחשבת והיא הפקה סדרת אודיו מסרטון וידאו נתון. מתברר שמדובר במשימה לא טריוויאלית והמאמר מציע גישה אלגנטית ודי פשוטה לבעייה זו. " + $br + $br + `
    "אז היום ב- #shorthebrewpapereviews סוקרים מאמר המציע שיטה להפקה אודיו מוידאו. בהינתן דאטהסט המכיל זוגות של סרטוני וידאו או אודיו המתאים המחברים משתמשים באנקודר מאומן של אודיו EnCodec שהופך את האודיו לייצוגו הלטנטי. מה זה ייצוג לטנטי של אות אודיו? " + $br + $br + `
    "למעשה זו סדרה של וקטורים שכל אחד מהם הוא השיכון (embedding) של מקטע (בזמן) של האות. בנוסף יש ל-EnCodec דקודר שמשחזר את האות מהייצוג הלטנטי שלו. המאמר גם משתמש במודלים שמטרתם להפיק ייצוג של וידאו (של כל פריים) כמו CLIP, ImageBind ו- ViT. " + $br + $br + `
    "אז מה בעצם הארכיטקטורה של FoleyGen ואיך מאמנים את המודל הזה? לכל זוג של וידאו ואודיו מעבירים את האודיו דרך האנדקור של EnCodec ואת הוידאו דרך האנקודר של דאטה ויזואלי (נגיד CLIP). כלומר כאן אודיו וידאו מויצגים באמצעות סדרה של וקטורי הייצוג של ״הטוקנים״ שמרכיבים אותם (פריים לוידאו ומקטע זמן לאודיו)." + $br + $br + `
    "לאחר מכן מאמנים טרנספורמר (מורכב מדקודר בלבד) שמטרתו לשחזר את ייצוג הטוקן הבא של אות אודיו בהינתן ייצוגי הטוקנים (של אודיו) הקודמים וייצוגי של טוקני הווידאו. הם בחנו כמה אופציות לגבי טוקנים של הווידאו של הטרנספורמר יכול לגשת: כל הטוקנים, רק הטוקנים שבאו לפני הזמן או את טוקני הווידאו הסמוכים בזמן. וזה וזה – פשוט ואלגנטי."
$d.Range($s, $e).Text = $newText2

# --- 4. Merge paragraphs "מודל של רכב..." + "נגלה היום.../אך האם..." into a
#        single paragraph whose only content is the new HuggingFace link.
$pModel = $d.Paragraphs.Item(8)
$pNigle = $d.Paragraphs.Item(9)
$markStart2 = $pModel.Range.End - 1
$markEnd2 = $pModel.Range.End
$d.Range($markStart2, $markEnd2).Delete()

$pMerged1 = $d.Paragraphs.Item(8)
$s2 = $pMerged1.Range.Start
$e2 = $pMerged1.Range.End - 1
$d.Range($s2, $e2).Text = "https://huggingface.co/papers/2309.10537"

# --- 5. Delete the three now-obsolete paragraphs right after the blank
#        paragraph that follows the paper link: the old arXiv abstract link,
#        the "נכתבה על ידי" credit line, and the blank line under it.
$pArxiv = $d.Paragraphs.Item(4)
$pAfter = $d.Paragraphs.Item(7)
$r3 = $d.Range($pArxiv.Range.Start, $pAfter.Range.Start)
$r3.Delete()
